# Pedidos.xlsx update - "atualizar Pedidos" / "Add files via upload"
# Fills in the previously-empty rows 142-156 (columns A:C) with new
# remessa/material/quantidade data, and updates the visible selection /
# scroll position accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to populate: Row, Remessa (col A), Material (col B), Quantidade (col C)
$rows = @(
    @(142, "80266429", "10247-ARI-I", 1),
    @(143, "80266430", "10493-ARI-I", 1),
    @(144, "80266430", "10255-ARI-I", 1),
    @(145, "80266431", "10256-ARI-I", 1),
    @(146, "80266431", "10636-ARI-I", 1),
    @(147, "80266434", "23359-ATE-I", 1),
    @(148, "80266435", "10527-ARI-I", 1),
    @(149, "80266436", "10030-XIN-I", 300),
    @(150, "80266437", "33720-SYN-I", 1),
    @(151, "80266438", "10547-ARI-I", 1),
    @(152, "80266438", "10541-ARI-I", 1),
    @(153, "80266439", "10548-ARI-I", 1),
    @(154, "80266440", "10060-ARI-I", 1),
    @(155, "80266440", "10040-ARI-I", 1),
    @(156, "80266441", "10548-ARI-I", 1)
)

# Remessa numbers (column A) look like plain numbers ("80266429"), but in
# the source data they are text. Assigning such a string straight to
# .Value would make Excel auto-convert it to a number, so instead we
# stage a text-typed FORMULA result ("="80266429"") in a scratch cell far
# away from the data, then copy only the *value* (PasteSpecial values)
# onto the real destination cell. A formula's text result carries its
# string-ness over via PasteSpecial without touching the destination
# cell's existing style/format (unlike NumberFormat="@", which would mint
# a brand-new cell style).
$scratch = $ws.Cells.Item(300, 10)

foreach ($r in $rows) {
    $rowNum = $r[0]

    $escapedRemessa = $r[1] -replace '"', '""'
    $scratch.Formula = '="' + $escapedRemessa + '"'
    $scratch.Copy()
    $ws.Cells.Item($rowNum, 1).PasteSpecial(-4163) | Out-Null
    $scratch.Clear()

    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}

$excel.CutCopyMode = 0

# Update view: selection now covers through row 156, and the window has
# scrolled down so row 112 is the top-left visible cell.
$ws.Range("A2:C156").Select()
$excel.ActiveWindow.ScrollRow = 112
